{"js": "// Replace the date line and every \"A\u00d7B=C\" multiplication answer in the\n// table with the new values, per the commit's regenerated-output diff.\nconst replacements = [\n  [\"2025-09-04 Thursday\", \"2025-09-05 Friday\"],\n  [\"810\u00d73=2430\", \"623\u00d73=1869\"],\n  [\"649\u00d74=2596\", \"127\u00d73=381\"],\n  [\"278\u00d76=1668\", \"924\u00d75=4620\"],\n  [\"650\u00d75=3250\", \"554\u00d74=2216\"],\n  [\"413\u00d76=2478\", \"112\u00d78=896\"],\n  [\"857\u00d77=5999\", \"627\u00d77=4389\"],\n  [\"257\u00d72=514\", \"194\u00d77=1358\"],\n  [\"422\u00d79=3798\", \"357\u00d75=1785\"],\n  [\"166\u00d73=498\", \"388\u00d73=1164\"],\n  [\"171\u00d73=513\", \"329\u00d75=1645\"],\n  [\"302\u00d77=2114\", \"777\u00d78=6216\"],\n  [\"810\u00d76=4860\", \"801\u00d77=5607\"],\n  [\"751\u00d74=3004\", \"425\u00d75=2125\"],\n  [\"995\u00d77=6965\", \"643\u00d77=4501\"],\n  [\"533\u00d78=4264\", \"667\u00d76=4002\"],\n  [\"106\u00d75=530\", \"821\u00d72=1642\"],\n  [\"413\u00d78=3304\", \"734\u00d77=5138\"],\n  [\"455\u00d75=2275\", \"521\u00d73=1563\"],\n  [\"346\u00d75=1730\", \"438\u00d74=1752\"],\n  [\"767\u00d73=2301\", \"288\u00d79=2592\"],\n  [\"492\u00d77=3444\", \"666\u00d77=4662\"],\n  [\"227\u00d74=908\", \"919\u00d75=4595\"],\n  [\"610\u00d76=3660\", \"766\u00d75=3830\"],\n  [\"906\u00d74=3624\", \"754\u00d75=3770\"],\n  [\"905\u00d73=2715\", \"361\u00d73=1083\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(\"Could not find text: \" + oldText);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every \"A\u00d7B=C\" multiplication answer in the\n# table with the new values, per the commit's regenerated-output diff.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2025-09-04 Thursday\", \"2025-09-05 Friday\"),\n  @(\"810\u00d73=2430\", \"623\u00d73=1869\"),\n  @(\"649\u00d74=2596\", \"127\u00d73=381\"),\n  @(\"278\u00d76=1668\", \"924\u00d75=4620\"),\n  @(\"650\u00d75=3250\", \"554\u00d74=2216\"),\n  @(\"413\u00d76=2478\", \"112\u00d78=896\"),\n  @(\"857\u00d77=5999\", \"627\u00d77=4389\"),\n  @(\"257\u00d72=514\", \"194\u00d77=1358\"),\n  @(\"422\u00d79=3798\", \"357\u00d75=1785\"),\n  @(\"166\u00d73=498\", \"388\u00d73=1164\"),\n  @(\"171\u00d73=513\", \"329\u00d75=1645\"),\n  @(\"302\u00d77=2114\", \"777\u00d78=6216\"),\n  @(\"810\u00d76=4860\", \"801\u00d77=5607\"),\n  @(\"751\u00d74=3004\", \"425\u00d75=2125\"),\n  @(\"995\u00d77=6965\", \"643\u00d77=4501\"),\n  @(\"533\u00d78=4264\", \"667\u00d76=4002\"),\n  @(\"106\u00d75=530\", \"821\u00d72=1642\"),\n  @(\"413\u00d78=3304\", \"734\u00d77=5138\"),\n  @(\"455\u00d75=2275\", \"521\u00d73=1563\"),\n  @(\"346\u00d75=1730\", \"438\u00d74=1752\"),\n  @(\"767\u00d73=2301\", \"288\u00d79=2592\"),\n  @(\"492\u00d77=3444\", \"666\u00d77=4662\"),\n  @(\"227\u00d74=908\", \"919\u00d75=4595\"),\n  @(\"610\u00d76=3660\", \"766\u00d75=3830\"),\n  @(\"906\u00d74=3624\", \"754\u00d75=3770\"),\n  @(\"905\u00d73=2715\", \"361\u00d73=1083\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
